$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Kurswahl" marks for rows 8-12 (columns E:O): the old scheme
# marked chosen subjects with "X"; the new scheme clears those and instead
# marks two (different) columns per row with "/". ---
# Row 8
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = "/"
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = "/"
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = ""
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = ""
$ws.Range("O8").Value = "/"
# Row 9
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = "/"
$ws.Range("G9").Value = ""
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = ""
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = ""
$ws.Range("L9").Value = "/"
$ws.Range("M9").Value = ""
$ws.Range("N9").Value = ""
$ws.Range("O9").Value = "/"
# Row 10
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = "/"
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = "/"
$ws.Range("K10").Value = ""
$ws.Range("L10").Value = ""
$ws.Range("M10").Value = "/"
$ws.Range("N10").Value = ""
$ws.Range("O10").Value = ""
# Row 11
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = "/"
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = "/"
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = "/"
$ws.Range("O11").Value = ""
# Row 12
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = "/"
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("L12").Value = "/"
$ws.Range("M12").Value = "/"
$ws.Range("N12").Value = ""
$ws.Range("O12").Value = ""
# --- Update the data validation on E8:Z49 from a fixed "X" list to a
# text-length check allowing 0-5 characters. ---
$validationRange = $ws.Range("E8:Z49")
$validationRange.Validation.Delete()
$validationRange.Validation.Add(6, 1, 1, 0, 5)
$validationRange.Validation.IgnoreBlank = $false
$validationRange.Validation.InCellDropdown = $true
$validationRange.Validation.ShowInput = $true
$validationRange.Validation.ShowError = $true

# --- Move the active selection in the frozen bottom-right pane from B8 to E8. ---
[void]$ws.Range("E8").Select()
